$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "292.12"
Set-TextValue $ws.Range("E2") "-5.41%"
Set-TextValue $ws.Range("D3") "39.93"
Set-TextValue $ws.Range("E3") "-3.07%"
Set-TextValue $ws.Range("D4") "5.034"
Set-TextValue $ws.Range("E4") "-2.90%"
Set-TextValue $ws.Range("D5") "0.07377"
Set-TextValue $ws.Range("E5") "-4.06%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D6") "4.309"
Set-TextValue $ws.Range("E6") "-0.19%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D7") "1.517"
Set-TextValue $ws.Range("E7") "-8.22%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.9222"
Set-TextValue $ws.Range("E8") "0.87%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Range("D9") "2.399"
Set-TextValue $ws.Range("E9") "-1.32%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D10") "0.1159"
Set-TextValue $ws.Range("E10") "-6.68%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D11") "0.1748"
Set-TextValue $ws.Range("E11") "-4.04%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.08683"
Set-TextValue $ws.Range("E12") "-5.62%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.04176"
Set-TextValue $ws.Range("E13") "-0.99%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.1052"
Set-TextValue $ws.Range("E14") "0.06%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001251"
Set-TextValue $ws.Range("E15") "-0.41%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D16") "0.006008"
Set-TextValue $ws.Range("E16") "4.41%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D17") "3.356"
Set-TextValue $ws.Range("E17") "0.31%"
Set-TextValue $ws.Range("E18") "-0.61%"
Set-TextValue $ws.Range("D19") "7.585"
Set-TextValue $ws.Range("E19") "2.76%"
Set-TextValue $ws.Range("D20") "0.1358"
Set-TextValue $ws.Range("E20") "-1.70%"
Set-TextValue $ws.Range("D21") "0.2817"
Set-TextValue $ws.Range("E21") "-0.07%"
Set-TextValue $ws.Range("D22") "0.03837"
Set-TextValue $ws.Range("E22") "-4.68%"
Set-TextValue $ws.Range("D23") "0.001289"
Set-TextValue $ws.Range("E23") "1.69%"
Set-TextValue $ws.Range("D24") "0.003607"
Set-TextValue $ws.Range("E24") "-11.82%"
Set-TextValue $ws.Range("E25") "0.58%"
Set-TextValue $ws.Range("D26") "0.0003744"
Set-TextValue $ws.Range("D38") "0.02317"
Set-TextValue $ws.Range("E38") "-9.49%"
Set-TextValue $ws.Range("D39") "0.04996"
Set-TextValue $ws.Range("E39") "-6.09%"
Set-TextValue $ws.Range("D40") "0.007733"
Set-TextValue $ws.Range("E40") "-1.55%"
Set-TextValue $ws.Range("D41") "0.1273"
Set-TextValue $ws.Range("E41") "-3.17%"
Set-TextValue $ws.Range("E42") "113.59%"
Set-TextValue $ws.Range("D43") "0.007434"
Set-TextValue $ws.Range("E43") "11.52%"
Set-TextValue $ws.Range("D44") "0.007924"
Set-TextValue $ws.Range("E44") "-1.09%"
Set-TextValue $ws.Range("D45") "0.3175"
Set-TextValue $ws.Range("E45") "3.36%"
Set-TextValue $ws.Range("D46") "0.00006491"
Set-TextValue $ws.Range("E46") "-3.49%"
Set-TextValue $ws.Range("D47") "0.00000000754"
Set-TextValue $ws.Range("E47") "0.54%"
Set-TextValue $ws.Range("E48") "43.39%"
Set-TextValue $ws.Range("D49") "0.004224"
Set-TextValue $ws.Range("E49") "36.11%"
Set-TextValue $ws.Range("D50") "0.00002112"
Set-TextValue $ws.Range("E50") "0.54%"
Set-TextValue $ws.Range("D51") "0.0002012"
Set-TextValue $ws.Range("E51") "0.54%"
